$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<other>"
$ws.Range("C2").Value = 51

# Row 3
$ws.Range("C3").Value = 58

# Row 4
$ws.Range("B4").Value = "<five>"
$ws.Range("C4").Value = 51

# Row 5
$ws.Range("B5").Value = "<is>"
$ws.Range("C5").Value = 51

# Row 6
$ws.Range("B6").Value = "<use>"
$ws.Range("C6").Value = 51

# Row 7
$ws.Range("C7").Value = 50

# Row 8
$ws.Range("B8").Value = "<make>"
$ws.Range("C8").Value = 52

# Row 9
$ws.Range("C9").Value = 52

# Row 10
$ws.Range("B10").Value = "<the>"
$ws.Range("C10").Value = 50

# Row 11
$ws.Range("B11").Value = "<an>"
$ws.Range("C11").Value = 55

# Row 12
$ws.Range("B12").Value = "<the>"
$ws.Range("C12").Value = 56

# Row 13
$ws.Range("B13").Value = "<enter>"
$ws.Range("C13").Value = 45

# Row 14
$ws.Range("B14").Value = "<by>"
$ws.Range("C14").Value = 49

# Row 15
$ws.Range("C15").Value = 54

# Row 16
$ws.Range("B16").Value = "<will>"
$ws.Range("C16").Value = 54

# Row 17
$ws.Range("B17").Value = "<come>"
$ws.Range("C17").Value = 55

# Row 18
$ws.Range("C18").Value = 44
